$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New students added to the attendance sheet (rows 19-21) ---
# (order matters so new shared strings land in the same position as the
#  target workbook: Rat Adrian, Bordas Norbert, Ivan Let Raul, Mailot Dorian)
$ws.Range("A19").Value = "Rat Adrian"
$ws.Range("A20").Value = "Bordas Norbert"

# Row 2's student name was reordered ("Let Raul Ivan" -> "Ivan Let Raul")
$ws.Range("A2").Value = "Ivan Let Raul"

$ws.Range("A21").Value = "Mailot Dorian"

# Copy the existing name-cell formatting (fill + full thin border) from the
# last pre-existing row down onto the 3 new rows.
$ws.Range("A18").Copy()
$ws.Range("A19:A21").PasteSpecial(-4122)

# New "sapt" attendance checkmarks (column C) for existing students.
$ws.Range("C3").Value = $true
$ws.Range("C6").Value = $true
$ws.Range("C7").Value = $true
$ws.Range("C11").Value = $true

# Attendance checkmarks for the 3 new students.
$ws.Range("C19").Value = $true
$ws.Range("C20").Value = $true
$ws.Range("C21").Value = $true

# The middle new row (A20) loses its top border since it sits flush under
# row 19's bottom border.
$ws.Range("A20").Borders.Item(8).LineStyle = -4142

[void]$ws.Range("F18").Select()
